$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text-only interpretation for Price/Volume columns so that
# values like "1.000" or "0.07814" are not reinterpreted as numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "29.408.33"
$ws.Range("E2").Value = "  +0.18%  "
$ws.Range("D3").Value = "1.868.48"
$ws.Range("E3").Value = "  -0.50%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("E5").Value = "  +0.29%  "
$ws.Range("D6").Value = "0.7039"
$ws.Range("E6").Value = "  -2.39%  "
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("D8").Value = "0.07969"
$ws.Range("E8").Value = "  -0.60%  "
$ws.Range("D9").Value = "0.3134"
$ws.Range("E9").Value = "  -0.56%  "
$ws.Range("D10").Value = "24.46"
$ws.Range("D11").Value = "0.07814"
$ws.Range("E11").Value = "  -4.85%  "
$ws.Range("D12").Value = "1.917.02"
$ws.Range("E12").Value = "  +2.04%  "
$ws.Range("D13").Value = "93.74"
$ws.Range("E13").Value = "  -1.06%  "
$ws.Range("D14").Value = "5.171"
$ws.Range("E14").Value = "  -1.17%  "
$ws.Range("D15").Value = "0.7024"
$ws.Range("E15").Value = "  -1.51%  "
$ws.Range("D16").Value = "6.482"
$ws.Range("E16").Value = "  +0.95%  "
$ws.Range("D17").Value = "0.000008638"
$ws.Range("E17").Value = "  +1.60%  "
$ws.Range("D18").Value = "29.502.52"
$ws.Range("E18").Value = "  +0.52%  "
$ws.Range("D19").Value = "251.95"
$ws.Range("E19").Value = "  +3.49%  "
$ws.Range("D20").Value = "2.155.73"
$ws.Range("E20").Value = "  +1.32%  "
$ws.Range("D21").Value = "13.12"
$ws.Range("E21").Value = "  -1.22%  "
$ws.Range("D22").Value = "1.000"
$ws.Range("E22").Value = "  -0.14%  "
$ws.Range("D23").Value = "7.664"
$ws.Range("E23").Value = "  -1.25%  "
$ws.Range("D24").Value = "1.001"
$ws.Range("E24").Value = "  -0.08%  "
$ws.Range("D25").Value = "0.1552"
$ws.Range("E25").Value = "  -2.57%  "
$ws.Range("D26").Value = "9.003"
$ws.Range("E26").Value = "  -0.45%  "
$ws.Range("D27").Value = "161.62"
$ws.Range("E27").Value = "  -0.54%  "
$ws.Range("D28").Value = "18.79"
$ws.Range("E28").Value = "  +1.43%  "
$ws.Range("D29").Value = "1.507"
$ws.Range("E29").Value = "  +0.38%  "
$ws.Range("D30").Value = "4.311"
$ws.Range("E30").Value = "  -2.31%  "
$ws.Range("D31").Value = "4.261"
$ws.Range("E31").Value = "  -1.12%  "
$ws.Range("D32").Value = "1.212"
$ws.Range("E32").Value = "  +0.87%  "
$ws.Range("D34").Value = "1.907"
$ws.Range("D35").Value = "0.7579"
$ws.Range("E35").Value = "  -0.63%  "
$ws.Range("D36").Value = "1.182"
$ws.Range("E36").Value = "  +0.32%  "
$ws.Range("D37").Value = "2.709"
$ws.Range("E37").Value = "  +0.01%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "0.01879"
$ws.Range("E38").Value = "  +0.20%  "
$ws.Range("B39").Value = "Maker"
$ws.Range("C39").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D39").Value = "1.279.12"
$ws.Range("E39").Value = "  +0.10%  "
$ws.Range("D40").Value = "2.771"
$ws.Range("E40").Value = "  +0.71%  "
$ws.Range("D41").Value = "0.8950"
$ws.Range("E41").Value = "  -1.74%  "
$ws.Range("D42").Value = "109.83"
$ws.Range("D43").Value = "6.027"
$ws.Range("E43").Value = "  -6.49%  "
$ws.Range("D44").Value = "70.99"
$ws.Range("E44").Value = "  -4.32%  "
$ws.Range("E45").Value = "  -0.15%  "
$ws.Range("D46").Value = "2.042.35"
$ws.Range("E46").Value = "  +0.79%  "
$ws.Range("E47").Value = "  -3.46%  "
$ws.Range("D48").Value = "1.808"
$ws.Range("E48").Value = "  +0.60%  "
$ws.Range("D49").Value = "9.608"
$ws.Range("E49").Value = "  +1.07%  "
$ws.Range("D50").Value = "0.5178"
$ws.Range("E50").Value = "  -0.96%  "
$ws.Range("D51").Value = "0.4294"
$ws.Range("E51").Value = "  -1.25%  "
